$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.363.55'
$ws.Range('E2').Value = '  +1.97%  '
$ws.Range('D3').Value = '1.826.53'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  -0.05%  '
$ws.Range('D5').Value = "'314.23"
$ws.Range('E5').Value = '  +1.26%  '
$ws.Range('D6').Value = "'0.9999"
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').Value = "'0.4698"
$ws.Range('E7').Value = '  +5.37%  '
$ws.Range('D8').Value = "'0.3802"
$ws.Range('E8').Value = '  +3.60%  '
$ws.Range('E9').Value = '  +1.17%  '
$ws.Range('D10').Value = "'0.8762"
$ws.Range('E10').Value = '  +2.28%  '
$ws.Range('D11').Value = "'20.81"
$ws.Range('E11').Value = '  +0.83%  '
$ws.Range('D12').Value = '1.827.63'
$ws.Range('E12').Value = '  -3.05%  '
$ws.Range('D13').Value = "'6.700"
$ws.Range('E13').Value = '  +1.48%  '
$ws.Range('D14').Value = "'5.436"
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').Value = "'93.21"
$ws.Range('E15').Value = '  +0.83%  '
$ws.Range('D16').Value = "'0.07093"
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('E17').Value = '  -0.03%  '
$ws.Range('D18').Value = "'0.000008796"
$ws.Range('E18').Value = '  +0.75%  '
$ws.Range('E19').Value = '  -0.02%  '
$ws.Range('D20').Value = "'15.04"
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').Value = '27.350.65'
$ws.Range('E21').Value = '  +1.80%  '
$ws.Range('D22').Value = "'5.330"
$ws.Range('E22').Value = '  +3.47%  '
$ws.Range('D23').Value = "'10.96"
$ws.Range('E23').Value = '  +1.52%  '
$ws.Range('D24').Value = '2.054.36'
$ws.Range('E24').Value = '  -2.96%  '
$ws.Range('D25').Value = "'1.941"
$ws.Range('E25').Value = '  -2.52%  '
$ws.Range('E26').Value = '  -0.55%  '
$ws.Range('D27').Value = "'2.255"
$ws.Range('E27').Value = '  +3.66%  '
$ws.Range('D28').Value = "'18.62"
$ws.Range('E28').Value = '  +1.06%  '
$ws.Range('D29').Value = "'5.337"
$ws.Range('E29').Value = '  +2.58%  '
$ws.Range('D30').Value = "'117.24"
$ws.Range('E30').Value = '  +0.66%  '
$ws.Range('D31').Value = "'0.08987"
$ws.Range('E31').Value = '  +2.06%  '
$ws.Range('D32').Value = "'0.7919"
$ws.Range('E32').Value = '  +5.93%  '
$ws.Range('D33').Value = "'1.193"
$ws.Range('E33').Value = '  +2.08%  '
$ws.Range('D34').Value = "'4.544"
$ws.Range('E34').Value = '  +1.78%  '
$ws.Range('D35').Value = "'2.936"
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  -0.04%  '
$ws.Range('D37').Value = "'1.102"
$ws.Range('E37').Value = '  +1.44%  '
$ws.Range('E38').Value = '  +0.69%  '
$ws.Range('D39').Value = "'0.05253"
$ws.Range('E39').Value = '  +1.08%  '
$ws.Range('D40').Value = "'7.301"
$ws.Range('E40').Value = '  +3.85%  '
$ws.Range('D41').Value = "'0.5336"
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('D42').Value = "'2.370"
$ws.Range('E42').Value = '  +20.58%  '
$ws.Range('D43').Value = "'2.897"
$ws.Range('E43').Value = '  +1.29%  '
$ws.Range('D44').Value = "'0.1705"
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('D45').Value = "'8.646"
$ws.Range('E45').Value = '  +2.22%  '
$ws.Range('D46').Value = "'0.5100"
$ws.Range('E46').Value = '  -0.85%  '
$ws.Range('D47').Value = "'10.63"
$ws.Range('E47').Value = '  +1.13%  '
$ws.Range('D48').Value = "'105.58"
$ws.Range('E48').Value = '  -0.05%  '
$ws.Range('D49').Value = "'1.682"
$ws.Range('E49').Value = '  +0.91%  '
$ws.Range('D50').Value = "'0.9995"
$ws.Range('E50').Value = '  +0.01%  '
$ws.Range('D51').Value = "'0.06386"
